$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.569.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.19%  "

$ws.Range("D3").Value = "'2.976.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.60%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'379.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.58%  "

$ws.Range("D6").Value = "'104.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.87%  "

$ws.Range("D7").Value = "'0.541"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.67%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "'0.594"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.11%  "

$ws.Range("D10").Value = "'37.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.50%  "

$ws.Range("E11").Value = "  +0.26%  "

$ws.Range("D12").Value = "'0.0845"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.07%  "

$ws.Range("D13").Value = "'3.448.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.94%  "

$ws.Range("D14").Value = "'18.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.39%  "

$ws.Range("D15").Value = "'7.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.39%  "

$ws.Range("D16").Value = "'2.979.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.87%  "

$ws.Range("D17").Value = "'0.973"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.70%  "

$ws.Range("D18").Value = "'51.549.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.33%  "

$ws.Range("D19").Value = "'3.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.56%  "

$ws.Range("E20").Value = "  +4.36%  "

$ws.Range("D21").Value = "'12.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.12%  "

$ws.Range("D22").Value = "'0.0₃0963"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.69%  "

$ws.Range("D23").Value = "'69.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.03%  "

$ws.Range("D24").Value = "'262.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.83%  "

$ws.Range("E25").Value = "  +7.36%  "

$ws.Range("D26").Value = "'8.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +19.08%  "

$ws.Range("D27").Value = "'7.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +24.53%  "

$ws.Range("D28").Value = "'0.117"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +15.32%  "

$ws.Range("E29").Value = "  -1.98%  "

$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").Value = "'25.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.66%  "

$ws.Range("D32").Value = "'9.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("D33").Value = "'35.06"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.50%  "

$ws.Range("D34").Value = "'51.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.56%  "

$ws.Range("E35").Value = "  -1.74%  "

$ws.Range("D36").Value = "'0.0449"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.97%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("E38").Value = "  +2.52%  "

$ws.Range("D39").Value = "'17.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.47%  "

$ws.Range("D40").Value = "'2.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.02%  "

$ws.Range("D41").Value = "'1.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.06%  "

$ws.Range("D42").Value = "'0.116"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.74%  "

$ws.Range("D43").Value = "'126.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.80%  "

$ws.Range("D44").Value = "'21.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.57%  "

$ws.Range("D45").Value = "'0.283"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +19.94%  "

$ws.Range("D46").Value = "'2.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.82%  "

$ws.Range("E47").Value = "  +3.24%  "

$ws.Range("D48").Value = "'2.036.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.37%  "

$ws.Range("E49").Value = "  +3.52%  "

$ws.Range("D50").Value = "'0.0333"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.21%  "

$ws.Range("D51").Value = "'58.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.78%  "
